# Replace the date heading and all 100 math-problem cell texts.
# Every "old" value below is unique in the document, and no "new"
# value collides with another pair's "old" value, so a sequence of
# document-wide Find/Replace (wdReplaceAll) calls is unambiguous.
$d = $word.ActiveDocument

$replacements = @(
    @("2023-07-30 Sunday", "2023-07-31 Monday"),
    @("61+19=", "91-31="),
    @("12+37=", "94-81="),
    @("20+73=", "18+9="),
    @("31-25=", "91-59="),
    @("59-18=", "32-27="),
    @("10+89=", "6+42="),
    @("3+31=", "38-28="),
    @("72-13=", "9+15="),
    @("44+25=", "64-62="),
    @("88-73=", "91+6="),
    @("92-65=", "89-83="),
    @("80+8=", "69+26="),
    @("94-83=", "78-10="),
    @("91-43=", "90-25="),
    @("69-68=", "4+37="),
    @("34-30=", "89-86="),
    @("41-36=", "10+49="),
    @("55+7=", "50+12="),
    @("94+0=", "69+10="),
    @("81+0=", "12+39="),
    @("25+26=", "73+13="),
    @("99-42=", "54-3="),
    @("52+25=", "33-4="),
    @("10+77=", "73-64="),
    @("41+51=", "44+2="),
    @("87-23=", "40+35="),
    @("20-2=", "19+22="),
    @("0+42=", "49-47="),
    @("14+33=", "90-60="),
    @("2+26=", "8+73="),
    @("65-46=", "44-18="),
    @("3+24=", "18+24="),
    @("98-23=", "84+8="),
    @("17+78=", "88-22="),
    @("85-3=", "6+52="),
    @("5+29=", "38-6="),
    @("9+65=", "39+44="),
    @("33-9=", "78-10="),
    @("48+40=", "72-10="),
    @("5+28=", "78-14="),
    @("69-43=", "40+26="),
    @("49+29=", "45+48="),
    @("71-38=", "65-14="),
    @("61-20=", "55+21="),
    @("63-44=", "50-40="),
    @("43-34=", "92-39="),
    @("10+76=", "93+3="),
    @("6+37=", "80+2="),
    @("3+81=", "53+12="),
    @("7+0=", "20-19="),
    @("30+62=", "27+46="),
    @("21-7=", "84-49="),
    @("87-39=", "39-39="),
    @("61-34=", "2+45="),
    @("35+54=", "93-38="),
    @("67-22=", "19+6="),
    @("77-42=", "89-54="),
    @("34+11=", "59-50="),
    @("32+42=", "16+31="),
    @("60-9=", "21+74="),
    @("80+4=", "9+51="),
    @("70-56=", "68-56="),
    @("70-44=", "46-39="),
    @("36+55=", "67+27="),
    @("47+51=", "19+76="),
    @("24+58=", "51+22="),
    @("23+29=", "49+35="),
    @("72-11=", "2+3="),
    @("14-0=", "66-63="),
    @("42-35=", "97-74="),
    @("18+54=", "23+0="),
    @("78+11=", "30-13="),
    @("86-53=", "37-2="),
    @("1+64=", "76-67="),
    @("46+30=", "41-16="),
    @("96-90=", "40-9="),
    @("10+71=", "75-10="),
    @("33+25=", "36-30="),
    @("59-27=", "41+42="),
    @("97-54=", "66+28="),
    @("14+49=", "0+87="),
    @("77+7=", "93-70="),
    @("95-74=", "87+2="),
    @("34-6=", "49+38="),
    @("48+0=", "49-43="),
    @("83-35=", "37+46="),
    @("30-4=", "33+63="),
    @("16+46=", "7+6="),
    @("84-68=", "87-4="),
    @("20-5=", "85-13="),
    @("37-19=", "52-44="),
    @("43-2=", "12-5="),
    @("68-22=", "42-32="),
    @("30+37=", "28+21="),
    @("95-47=", "25+33="),
    @("1+10=", "78-20="),
    @("60+30=", "22-9="),
    @("44+52=", "54+27="),
    @("23+62=", "81-7="),
    @("4+49=", "82-40="),
)

$count = 0
foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $findRange = $d.Content
    $found = $findRange.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if ($found) {
        $count++
    } else {
        Write-Host "NOT FOUND: $old"
    }
}

Write-Host "Total replacements: $count of $($replacements.Count)"
